# Apply the TestData.xlsx edits:
#  - Update a few password values in column B
#  - Replace the shared "12345678" string used by B8 with "wqerewr"
#  - Move the active selection from D8 to G4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

# Update password values for a few rows
$ws.Range("B3").Value = 12345
$ws.Range("B4").Value = 1378
$ws.Range("B5").Value = 12678

# B8 keeps its text-quoted style (quotePrefix) - use a leading apostrophe
# so Excel treats it as literal text and preserves the cell's format.
$ws.Range("B8").Value = "'wqerewr"

# Update the active selection on the sheet
$ws.Range("G4").Select()
